# Sample_Registrar_Masters.xlsx edit:
#  - Mark three header cells as "required" fields by appending an asterisk.
#  - Introduce a new (currently blank) column R between the last data column
#    (Q) and the decorative trailing columns (S:U), matching the formatting
#    already used by S1/S2.
#  - Re-tint the sheet's thin-box borders: the main data grid (A1:Q2) moves
#    to the slightly darker #888888, and the accent block (R1:U2) moves to
#    the even darker #757575.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R -----------------------------------------------------
# Clone the formatting of the neighbouring (currently last) decorative
# column so the new blank cells pick up the same style family (font/
# border) that S1/S2 already use, rather than reverting to a default style.
$ws.Range("S1").Copy()
$ws.Range("R1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("S2").Copy()
$ws.Range("R2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Header text tweaks (flag required fields with a trailing "*") ----
$ws.Range("A1").Value = "Registrar Name*"
$ws.Range("B1").Value = "SEBI Regn. ID*"
$ws.Range("Q1").Value = "Company Master Id*"

# --- Border recolouring -------------------------------------------------
# Main data grid: #9A9A9A -> #888888
$ws.Range("A1:Q2").Borders.Color = 8947848
# Accent block (now including the new column R): #888888 -> #757575
$ws.Range("R1:U2").Borders.Color = 7697781
